$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "B4"  = 19
    "B7"  = 16
    "B9"  = 22
    "B11" = 32
    "B12" = 40
    "B13" = 55
    "B14" = 47
    "B15" = 51
    "B16" = 60
    "B19" = 66
    "B23" = 80
    "B25" = 79
    "B27" = 76
    "B29" = 80
    "B30" = 65
    "B33" = 60
    "B34" = 64
    "B35" = 57
    "B38" = 46
    "B40" = 38
    "B41" = 27
    "B42" = 33
    "B46" = 31
    "B47" = 24
    "B50" = 25
    "B51" = 23
    "B52" = 20
    "B53" = 16
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
